# oox smartart, linear layout: correctly scale spacings wrt constraints and rules
#
# The diagram's graphic frame on slide 1 is moved further down the slide
# (its vertical offset grows from 1407600 EMU to 2847600 EMU, i.e. from
# 110.8346pt to 224.2205pt), while its width/height and horizontal
# position stay the same.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)

# 2847600 EMU / 914400 EMU-per-inch * 72 pt-per-inch = 224.2205 pt
$sh.Top = 2847600 / 914400 * 72
